# Apply "Update of numbers in paper" revisions to the oc_reg regression table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: "FE" row - YES -> checkmark (create "\checkmark" shared string first) ---
$ws.Range("B13").Value = "\checkmark"

# --- Header row (row 2): rename the two panel headers -----------------
# B2:E2 merged -> "Take-up (choice arms)"
$ws.Range("B2").Value = "Take-up (choice arms)"
# F2:G2 merged -> "Financing Cost (hte)"
$ws.Range("F2").Value = "Financing Cost (hte)"

# --- Row 5: "OC" coefficient label becomes a static "OC (dummy)" ------
# (was driven by an external-link formula; replace with plain text)
$ws.Range("A5").Value = "OC (dummy)"

# --- Row 13 label: "FE" -> "Branch/Day FE" -----------------------------
$ws.Range("A13").Value = "Branch/Day FE"

# finish filling in the rest of row 13 with checkmarks
$ws.Range("C13:G13").Value = "\checkmark"
# keep these centered (matches surrounding cell formatting)
$ws.Range("C13:G13").HorizontalAlignment = -4108

# --- Row 14: "Controls" row - NO -> blank, YES -> checkmark ------------
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("D14").Value = "\checkmark"
$ws.Range("E14").Value = "\checkmark"
$ws.Range("F14").ClearContents()
$ws.Range("G14").Value = "\checkmark"

$wb.Save()
